$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.091.23"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.901.87"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'361.53"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'104.45"
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("D7").Value = "'0.540"
$ws.Range("E7").Value = "  -3.91%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  -5.88%  "
$ws.Range("D10").Value = "'36.77"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "'0.0833"
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("D13").Value = "'18.43"
$ws.Range("E13").Value = "  -5.14%  "
$ws.Range("D14").Value = "3.357.26"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "'7.34"
$ws.Range("E15").Value = "  -4.49%  "
$ws.Range("D16").Value = "2.905.12"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "51.028.75"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "'3.30"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").Value = "'7.20"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").Value = "'12.99"
$ws.Range("E21").Value = "  -5.71%  "
$ws.Range("D22").Value = "0.0₃0943"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").Value = "'68.09"
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("D24").Value = "'258.08"
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").Value = "'2.68"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").Value = "'0.173"
$ws.Range("E26").Value = "  -6.10%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'25.85"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "'7.17"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'9.92"
$ws.Range("E32").Value = "  -5.02%  "
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").Value = "'34.80"
$ws.Range("E34").Value = "  -6.84%  "
$ws.Range("D35").Value = "'50.59"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "'0.0421"
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("D38").Value = "'2.81"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").Value = "'3.13"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("D40").Value = "'16.91"
$ws.Range("E40").Value = "  -6.90%  "
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  -6.55%  "
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("D43").Value = "'22.35"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "'118.77"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'2.13"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "2.068.55"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("E47").Value = "  -6.91%  "
$ws.Range("D48").Value = "'2.27"
$ws.Range("E48").Value = "  -8.35%  "
$ws.Range("D49").Value = "3.193.77"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  -6.35%  "
$ws.Range("D51").Value = "'0.0309"
$ws.Range("E51").Value = "  -8.42%  "
